# Updates cryptos list data (Price / Volume(1h) columns, plus a few
# re-ordered Coin/Link/Price rows) to match the refreshed GitHub Actions
# snapshot. Values that look like plain numbers (e.g. "234.36") are
# prefixed with a leading apostrophe so Excel stores them as text, just
# like the source workbook's existing cells (t="inlineStr").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.371.74'
$ws.Range("E2").Value = '  -0.97%  '
$ws.Range("D3").Value = '2.075.01'
$ws.Range("E3").Value = '  -0.53%  '
$ws.Range("E4").Value = '  +0.12%  '
$ws.Range("D5").Value = '''234.36'
$ws.Range("E5").Value = '  -1.18%  '
$ws.Range("D6").Value = '''0.624'
$ws.Range("E6").Value = '  +1.16%  '
$ws.Range("D8").Value = '''57.05'
$ws.Range("E8").Value = '  -2.72%  '
$ws.Range("E9").Value = '  -1.07%  '
$ws.Range("E10").Value = '  -0.18%  '
$ws.Range("E11").Value = '  +0.54%  '
$ws.Range("D12").Value = '2.379.62'
$ws.Range("E12").Value = '  -0.42%  '
$ws.Range("D13").Value = '''14.68'
$ws.Range("E13").Value = '  +0.75%  '
$ws.Range("D14").Value = '''20.78'
$ws.Range("E14").Value = '  -1.35%  '
$ws.Range("D15").Value = '''0.778'
$ws.Range("E15").Value = '  -0.63%  '
$ws.Range("E16").Value = '  -2.54%  '
$ws.Range("D17").Value = '2.075.21'
$ws.Range("D18").Value = '37.291.64'
$ws.Range("E18").Value = '  -0.85%  '
$ws.Range("E19").Value = '  +2.22%  '
$ws.Range("D20").Value = '''69.46'
$ws.Range("E20").Value = '  +1.03%  '
$ws.Range("D21").Value = '0.0₃0814'
$ws.Range("E21").Value = '  -0.29%  '
$ws.Range("D22").Value = '''226.81'
$ws.Range("E22").Value = '  +0.81%  '
$ws.Range("E23").Value = '  +0.05%  '
$ws.Range("D24").Value = '''2.43'
$ws.Range("E24").Value = '  -0.02%  '
$ws.Range("E25").Value = '  -3.17%  '
$ws.Range("D26").Value = '''167.18'
$ws.Range("E26").Value = '  +2.24%  '
$ws.Range("E27").Value = '  -1.13%  '
$ws.Range("E28").Value = '  +2.85%  '
$ws.Range("D29").Value = '''19.08'
$ws.Range("E29").Value = '  -1.86%  '
$ws.Range("E30").Value = '  -5.47%  '
$ws.Range("E31").Value = '  -1.13%  '
$ws.Range("E32").Value = '  -0.63%  '
$ws.Range("E33").Value = '  -2.30%  '
$ws.Range("D34").Value = '''4.57'
$ws.Range("E34").Value = '  +1.85%  '
$ws.Range("E35").Value = '  -4.35%  '
$ws.Range("E36").Value = '  +0.21%  '
$ws.Range("E37").Value = '  -0.40%  '
$ws.Range("E38").Value = '  -3.66%  '
$ws.Range("D39").Value = '''5.69'
$ws.Range("E39").Value = '  -4.76%  '
$ws.Range("E40").Value = '  -0.25%  '
$ws.Range("D41").Value = '''4.42'
$ws.Range("E41").Value = '  +2.94%  '
$ws.Range("D42").Value = '1.476.06'
$ws.Range("E42").Value = '  -0.47%  '
$ws.Range("D43").Value = '''0.0942'
$ws.Range("E43").Value = '  -2.93%  '
$ws.Range("D44").Value = '''96.43'
$ws.Range("E44").Value = '  +0.90%  '
$ws.Range("B45").Value = 'TrustWalletToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D45").Value = '''1.17'
$ws.Range("E45").Value = '  +2.60%  '
$ws.Range("B46").Value = 'VeChain'
$ws.Range("C46").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D46").Value = '''0.0212'
$ws.Range("E46").Value = '  -0.50%  '
$ws.Range("E47").Value = '  -0.97%  '
$ws.Range("B48").Value = 'InjectiveProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D48").Value = '''15.06'
$ws.Range("E48").Value = '  -8.68%  '
$ws.Range("B49").Value = 'FraxShare'
$ws.Range("C49").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D49").Value = '''7.16'
$ws.Range("E49").Value = '  -2.61%  '
$ws.Range("E50").Value = '  +0.83%  '
$ws.Range("D51").Value = '2.266.99'
$ws.Range("E51").Value = '  -0.50%  '
